$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.104.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.079.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.675'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.96'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +15.54%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.391'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '61.46'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0799'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.108'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.388.25'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.821'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.082.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.098.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.94%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '74.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0929'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.126'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.92%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0637'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0917'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.31'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.71%  '
$ws.Range('B39').Value = 'Cronos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.116'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +27.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.93%  '
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +39.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0227'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.20%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.82'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.97'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.309.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.69%  '
